$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7, column A: was stored as text "71652621" -> should become a real number
$ws.Cells.Item(7, 1).Value = 71652621

# New row 8: payment 76442781 (Cash) 2025-08-15T09:48:27
# Column A must stay a text value (looks numeric) like the existing A7 used to be,
# so force Text formatting while entering it, then clear the style back to Normal
# so no stray formatting is left behind on the cell.
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "76442781"
$ws.Cells.Item(8, 1).Style = "Normal"

$ws.Cells.Item(8, 2).Value = 20
$ws.Cells.Item(8, 3).Value = "Cash"
$ws.Cells.Item(8, 4).Value = "2025-08-15T09:48:27"
